# Update the "LOS Galacticos" roster sheet with refreshed player / position / team data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @(
    @(2,  "Tyus Jones",             "PG",       "Phoenix Suns"),
    @(3,  "Keyonte George",         "PG,SG",    "Utah Jazz"),
    @(4,  "Anthony Edwards",        "SG,SF",    "Minnesota Timberwolves"),
    @(5,  "Fred VanVleet",          "PG",       "Houston Rockets"),
    @(6,  "Amen Thompson",          "SG,SF",    "Houston Rockets"),
    @(7,  "P.J. Washington",        "PF",       "Dallas Mavericks"),
    @(8,  "Jayson Tatum",           "SF,PF",    "Boston Celtics"),
    @(9,  "Jaren Jackson Jr.",      "PF,C",     "Memphis Grizzlies"),
    @(10, "Zion Williamson",        "PF,C",     "New Orleans Pelicans"),
    @(11, "Giannis Antetokounmpo",  "PF,C",     "Milwaukee Bucks"),
    @(12, "Bradley Beal",           "PG,SG,SF", "Phoenix Suns"),
    @(13, "James Harden",           "PG,SG",    "LA Clippers"),
    @(14, "Anfernee Simons",        "PG,SG",    "Portland Trail Blazers"),
    @(15, "Paul George",            "SG,SF,PF", "Philadelphia 76ers"),
    @(16, "Jaylin Williams",        "PF,C",     "Oklahoma City Thunder"),
    @(17, "Jonathan Kuminga",       "SF,PF",    "Golden State Warriors"),
    @(18, "Goga Bitadze",           "C",        "Orlando Magic"),
    @(19, "Ivica Zubac",            "C",        "LA Clippers")
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
}
